# Projeto v0.9 - Criado executavel para testes
# Simplify the contact list down to a single row ("EU" / "999999999"),
# dropping the "Mensagem" column and the extra sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the sample data in row 2 with the new single contact.
$ws.Range("A2").Value = "EU"
$ws.Range("B2").Value = "999999999"

# Remove the now-unneeded rows 3 and 4 (Pessoa 2 / Pessoa 3).
$ws.Rows("3:4").Delete() | Out-Null

# Remove the "Mensagem" column entirely.
$ws.Columns("C:C").Delete() | Out-Null

# Match the saved selection from the authored workbook.
$ws.Range("B5").Select() | Out-Null
